$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dainandin Nond")
$ws.Range("H31").Value = 7
